# Updated cryptos list on Mon Aug 28 07:20:38 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.074.99"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.650.50"
$ws.Range("E3").Value = "  -0.84%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.39%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'217.64"
$ws.Range("E5").Value = "  -0.61%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.5211"
$ws.Range("E6").Value = "  -2.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.38%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2617"
$ws.Range("E8").Value = "  -1.19%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06288"
$ws.Range("E9").Value = "  -2.10%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'20.50"
$ws.Range("E10").Value = "  -0.48%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  -0.60%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'4.478"
$ws.Range("E12").Value = "  -2.02%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.598.76"
$ws.Range("E13").Value = "  -3.96%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'1.878.00"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.5542"
$ws.Range("E15").Value = "  +0.39%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "'0.0₅7993"
$ws.Range("E16").Value = "  -2.64%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'64.78"
$ws.Range("E17").Value = "  -1.34%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'26.081.41"
$ws.Range("E18").Value = "  -0.87%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.38%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'4.627"
$ws.Range("E20").Value = "  -1.41%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'194.31"
$ws.Range("E21").Value = "  +0.34%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'10.05"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "'5.944"
$ws.Range("E23").Value = "  -1.63%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  -0.36%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'146.54"
$ws.Range("E25").Value = "  +0.47%  "

# Row 26 - Stellar
$ws.Range("D26").Value = "'0.1203"
$ws.Range("E26").Value = "  -2.39%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'7.171"
$ws.Range("E27").Value = "  -0.40%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'15.88"
$ws.Range("E28").Value = "  -1.57%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'1.475"
$ws.Range("E29").Value = "  -0.48%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "'0.05612"
$ws.Range("E30").Value = "  -4.15%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.28%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'3.476"
$ws.Range("E32").Value = "  -4.04%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'3.348"
$ws.Range("E33").Value = "  +1.97%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.591"
$ws.Range("E34").Value = "  -1.09%  "

# Row 35 - was ARBITRUM, now MXToken
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").Value = "'2.800"
$ws.Range("E35").Value = "  -0.91%  "

# Row 36 - was MXToken, now ARBITRUM
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9498"
$ws.Range("E36").Value = "  -1.58%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "'2.410"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "'0.5648"
$ws.Range("E38").Value = "  -2.73%  "

# Row 39 - was FraxShare, now VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01583"
$ws.Range("E39").Value = "  -1.62%  "

# Row 40 - was VeChain, now FraxShare
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.951"
$ws.Range("E40").Value = "  +1.10%  "

# Row 41 - Maker
$ws.Range("D41").Value = "'1.056.65"
$ws.Range("E41").Value = "  +0.58%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'0.8410"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44 - Quant
$ws.Range("D44").Value = "'102.75"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "'1.789.47"
$ws.Range("E45").Value = "  -0.77%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'57.16"
$ws.Range("E46").Value = "  -0.95%  "

# Row 47 - BabyDogeCoin
$ws.Range("E47").Value = "  +0.54%  "

# Row 48 - Frax
$ws.Range("D48").Value = "'1.010"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.05299"
$ws.Range("E49").Value = "  +2.58%  "

# Row 50 - Mantle
$ws.Range("D50").Value = "'0.4339"
$ws.Range("E50").Value = "  -1.00%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "'7.946"
$ws.Range("E51").Value = "  -1.39%  "
